$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. New column R header "backup" - copy header format (bold/border/alignment) from Q1
$ws.Range("R1").Value = "backup"
$ws.Range("Q1").Copy()
$ws.Range("R1").PasteSpecial(-4122)

# 2. Fill R2:R401 with 0 (new "backup" column data for existing rows)
$ws.Range("R2:R401").Value = 0

# 3. A few rows got a non-zero backup value instead of the default 0
$ws.Range("R108").Value = 2
$ws.Range("R197").Value = 2
$ws.Range("R206").Value = 2

# 4. Several existing detect_structure (Q) values were recomputed to 0
$ws.Range("Q30").Value = 0
$ws.Range("Q42").Value = 0
$ws.Range("Q47").Value = 0
$ws.Range("Q49").Value = 0
$ws.Range("Q55").Value = 0

# 5. Last existing row's isPivot (O) flag was recomputed to 2
$ws.Range("O402").Value = 2

# 6. Last existing row also gets its new backup column value
$ws.Range("R402").Value = 0

# 7. Append 6 new monthly rows (403-408), carrying forward row 402's formatting
$ws.Range("A402:Q402").Copy()
$ws.Range("A403:Q408").PasteSpecial(-4122)

$ws.Cells.Item(403,1).Value = 45474
$ws.Cells.Item(403,2).Value = 989.75
$ws.Cells.Item(403,3).Value = 1179
$ws.Cells.Item(403,4).Value = 967.2000122070312
$ws.Cells.Item(403,5).Value = 1156.650024414062
$ws.Cells.Item(403,7).Value = 272457360
$ws.Cells.Item(403,8).Value = 2024
$ws.Cells.Item(403,9).Value = 7
$ws.Cells.Item(403,10).Value = 1
$ws.Cells.Item(403,11).Value = 0
$ws.Cells.Item(403,12).Value = 0
$ws.Cells.Item(403,13).Value = 0
$ws.Cells.Item(403,14).Value = 27
$ws.Cells.Item(403,15).Value = 1
$ws.Cells.Item(403,16).Value = 0
$ws.Cells.Item(403,17).Value = 2

$ws.Cells.Item(404,1).Value = 45505
$ws.Cells.Item(404,2).Value = 1167
$ws.Cells.Item(404,3).Value = 1176
$ws.Cells.Item(404,4).Value = 1008.400024414062
$ws.Cells.Item(404,5).Value = 1111.349975585938
$ws.Cells.Item(404,7).Value = 284783154
$ws.Cells.Item(404,8).Value = 2024
$ws.Cells.Item(404,9).Value = 8
$ws.Cells.Item(404,10).Value = 1
$ws.Cells.Item(404,11).Value = 0
$ws.Cells.Item(404,12).Value = 0
$ws.Cells.Item(404,13).Value = 0
$ws.Cells.Item(404,14).Value = 31
$ws.Cells.Item(404,15).Value = 0
$ws.Cells.Item(404,16).Value = 0
$ws.Cells.Item(404,17).Value = 0

$ws.Cells.Item(405,1).Value = 45536
$ws.Cells.Item(405,2).Value = 1105
$ws.Cells.Item(405,3).Value = 1105
$ws.Cells.Item(405,4).Value = 949.2000122070312
$ws.Cells.Item(405,5).Value = 974.6500244140625
$ws.Cells.Item(405,7).Value = 296458308
$ws.Cells.Item(405,8).Value = 2024
$ws.Cells.Item(405,9).Value = 9
$ws.Cells.Item(405,10).Value = 1
$ws.Cells.Item(405,11).Value = 0
$ws.Cells.Item(405,12).Value = 0
$ws.Cells.Item(405,13).Value = 0
$ws.Cells.Item(405,14).Value = 35
$ws.Cells.Item(405,15).Value = 0
$ws.Cells.Item(405,16).Value = 0
$ws.Cells.Item(405,17).Value = 0

$ws.Cells.Item(406,1).Value = 45566
$ws.Cells.Item(406,2).Value = 976.9000244140625
$ws.Cells.Item(406,3).Value = 984.5
$ws.Cells.Item(406,4).Value = 825.7000122070312
$ws.Cells.Item(406,5).Value = 834.0499877929688
$ws.Cells.Item(406,7).Value = 274830221
$ws.Cells.Item(406,8).Value = 2024
$ws.Cells.Item(406,9).Value = 10
$ws.Cells.Item(406,10).Value = 1
$ws.Cells.Item(406,11).Value = 0
$ws.Cells.Item(406,12).Value = 0
$ws.Cells.Item(406,13).Value = 0
$ws.Cells.Item(406,14).Value = 40
$ws.Cells.Item(406,15).Value = 0
$ws.Cells.Item(406,16).Value = 0
$ws.Cells.Item(406,17).Value = 1

$ws.Cells.Item(407,1).Value = 45597
$ws.Cells.Item(407,2).Value = 847.9500122070312
$ws.Cells.Item(407,3).Value = 847.9500122070312
$ws.Cells.Item(407,4).Value = 759.2000122070312
$ws.Cells.Item(407,5).Value = 786.4500122070312
$ws.Cells.Item(407,7).Value = 255661808
$ws.Cells.Item(407,8).Value = 2024
$ws.Cells.Item(407,9).Value = 11
$ws.Cells.Item(407,10).Value = 1
$ws.Cells.Item(407,11).Value = 0
$ws.Cells.Item(407,12).Value = 0
$ws.Cells.Item(407,13).Value = 0
$ws.Cells.Item(407,14).Value = 44
$ws.Cells.Item(407,15).Value = 0
$ws.Cells.Item(407,16).Value = 0
$ws.Cells.Item(407,17).Value = 2

$ws.Cells.Item(408,1).Value = 45627
$ws.Cells.Item(408,2).Value = 787.2999877929688
$ws.Cells.Item(408,3).Value = 820.3499755859375
$ws.Cells.Item(408,4).Value = 717.7000122070312
$ws.Cells.Item(408,5).Value = 733.6500244140625
$ws.Cells.Item(408,7).Value = 263012899
$ws.Cells.Item(408,8).Value = 2024
$ws.Cells.Item(408,9).Value = 12
$ws.Cells.Item(408,10).Value = 1
$ws.Cells.Item(408,11).Value = 0
$ws.Cells.Item(408,12).Value = 0
$ws.Cells.Item(408,13).Value = 0
$ws.Cells.Item(408,14).Value = 48
$ws.Cells.Item(408,15).Value = 0
$ws.Cells.Item(408,16).Value = 0
$ws.Cells.Item(408,17).Value = 0

Write-Host "edit complete"
